$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Order ID" test data (N2) is refreshed to a new, clean cart code.
# Force the cell to stay text (it holds a numeric-looking code like the
# original "603424"), then restore the default "Normal" style so the
# cell keeps the same formatting it started with.
$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "420565"
$ws.Range("N2").Style = "Normal"
